$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the numeric-looking text values in columns D and E stay as text
# (matching the workbook's original inlineStr storage) rather than being
# auto-converted to numbers by Excel.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '26.403.11'
$ws.Range('E2').Value = '  -0.72%  '
$ws.Range('D3').Value = '1.840.19'
$ws.Range('E3').Value = '  -0.98%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '261.43'
$ws.Range('E5').Value = '  -4.03%  '
$ws.Range('E6').Value = '  -0.04%  '
$ws.Range('D7').Value = '0.5183'
$ws.Range('E7').Value = '  -1.42%  '
$ws.Range('D8').Value = '0.3266'
$ws.Range('E8').Value = '  -3.54%  '
$ws.Range('E9').Value = '  -0.39%  '
$ws.Range('D10').Value = '18.58'
$ws.Range('E10').Value = '  -6.39%  '
$ws.Range('E11').Value = '  -2.33%  '
$ws.Range('D12').Value = '0.07759'
$ws.Range('E12').Value = '  +0.37%  '
$ws.Range('D13').Value = '1.822.35'
$ws.Range('E13').Value = '  -2.75%  '
$ws.Range('D14').Value = '87.51'
$ws.Range('E14').Value = '  -2.38%  '
$ws.Range('D15').Value = '4.996'
$ws.Range('E15').Value = '  -2.59%  '
$ws.Range('D16').Value = '0.9999'
$ws.Range('E16').Value = '  +0.01%  '
$ws.Range('E17').Value = '  -3.68%  '
$ws.Range('D18').Value = '0.9999'
$ws.Range('E18').Value = '  -0.15%  '
$ws.Range('D19').Value = '0.000007965'
$ws.Range('E19').Value = '  -0.18%  '
$ws.Range('D20').Value = '26.423.30'
$ws.Range('E20').Value = '  -0.77%  '
$ws.Range('D21').Value = '2.068.51'
$ws.Range('E21').Value = '  -3.01%  '
$ws.Range('D22').Value = '4.611'
$ws.Range('E22').Value = '  -2.46%  '
$ws.Range('D23').Value = '9.502'
$ws.Range('E23').Value = '  -4.55%  '
$ws.Range('D24').Value = '5.971'
$ws.Range('E24').Value = '  -2.31%  '
$ws.Range('D25').Value = '144.69'
$ws.Range('E25').Value = '  -0.60%  '
$ws.Range('D26').Value = '2.182'
$ws.Range('E26').Value = '  -7.71%  '
$ws.Range('D27').Value = '1.648'
$ws.Range('E27').Value = '  -0.36%  '
$ws.Range('E28').Value = '  -1.51%  '
$ws.Range('D29').Value = '111.65'
$ws.Range('E29').Value = '  -0.20%  '
$ws.Range('D30').Value = '4.167'
$ws.Range('E30').Value = '  -3.08%  '
$ws.Range('E31').Value = '  -4.32%  '
$ws.Range('D32').Value = '0.08693'
$ws.Range('E32').Value = '  -2.27%  '
$ws.Range('D33').Value = '0.04814'
$ws.Range('E33').Value = '  -2.11%  '
$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D34').Value = '0.7207'
$ws.Range('E34').Value = '  -0.76%  '
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').Value = '1.125'
$ws.Range('E35').Value = '  -2.55%  '
$ws.Range('D36').Value = '2.849'
$ws.Range('E36').Value = '  -1.08%  '
$ws.Range('E37').Value = '  -4.23%  '
$ws.Range('D38').Value = '0.01775'
$ws.Range('E38').Value = '  -3.79%  '
$ws.Range('D39').Value = '2.213'
$ws.Range('E39').Value = '  -5.15%  '
$ws.Range('D40').Value = '0.4802'
$ws.Range('E40').Value = '  -5.76%  '
$ws.Range('D41').Value = '0.9063'
$ws.Range('E41').Value = '  -2.75%  '
$ws.Range('D42').Value = '111.26'
$ws.Range('E42').Value = '  -4.08%  '
$ws.Range('D43').Value = '6.054'
$ws.Range('E43').Value = '  -1.49%  '
$ws.Range('E44').Value = '  -0.04%  '
$ws.Range('D45').Value = '7.690'
$ws.Range('E45').Value = '  -3.79%  '
$ws.Range('D46').Value = '0.05915'
$ws.Range('E46').Value = '  -0.35%  '
$ws.Range('D47').Value = '0.4141'
$ws.Range('E47').Value = '  -6.12%  '
$ws.Range('D48').Value = '8.995'
$ws.Range('E48').Value = '  -2.89%  '
$ws.Range('D49').Value = '34.94'
$ws.Range('E49').Value = '  -2.84%  '
$ws.Range('D50').Value = '0.1224'
$ws.Range('E50').Value = '  -7.53%  '
$ws.Range('D51').Value = '0.8839'
$ws.Range('E51').Value = '  +0.60%  '
